$d = $word.ActiveDocument

# --- Remove the old, mid-document "_GoBack" bookmark (it sat right after the
# "...artifact use roll d20..." run) -------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Insert a new "4" run (bold, en-US) immediately before the "Spades" run,
# then re-create the (now empty/collapsed) "_GoBack" bookmark right after it
# and right before "Spades" -------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Spades", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $insertionStart = $rng.Start

    $fourRange = $d.Range($insertionStart, $insertionStart)
    $fourRange.InsertBefore("4")

    $fourRunRange = $d.Range($insertionStart, $insertionStart + 1)
    $fourRunRange.Font.Bold = 1
    $fourRunRange.LanguageID = "en-US"

    $bookmarkPoint = $d.Range($insertionStart + 1, $insertionStart + 1)
    $d.Bookmarks.Add("_GoBack", $bookmarkPoint)
}
